$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the value in F3 (was 17.64) -> becomes blank/missing
$ws.Range("F3").ClearContents()

# 2. Delete the entire row for "RM 232" (row 26) - data row removed entirely
$ws.Rows.Item(26).Delete()

# 3. Delete the entire row for "SC 92" (originally row 28, now row 27
#    after the previous deletion shifted everything up by one)
$ws.Rows.Item(27).Delete()

# 4. After the two row deletions, a few previously-missing values in the
#    shifted rows are now populated (revealing underlying data):
#    - Row 26 ("SC 5"): E26 now has a value
$ws.Range("E26").Value = -5

#    - Row 27 ("SC 101"): E27 becomes missing (was -10)
$ws.Range("E27").ClearContents()

#    - Row 33 ("SC 232"): E33 and F33 now have values
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53
